$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the two new columns (I0, IF), matching the style
# (bold font, thin border, centered/top alignment) of the existing headers.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J78 (column I then column J)
$iVals = @(9,7,5,6,5,8,10,8,8,5,6,9,8,6,9,9,9,9,9,9,9,9,9,8,9,9,9,9,9,9,9,9,9,9,9,9,8,9,8,9,9,7,9,9,9,9,9,10,9,9,9,8,9,10,8,7,9,9,9,9,9,9,8,8,9,9,9,8,9,9,7,8,8,9,4,3,6)
$jVals = @(9,7,6,6,6,8,10,8,8,6,6,9,8,6,9,9,9,9,10,9,9,9,9,8,9,9,9,9,9,9,9,9,9,9,10,9,9,9,8,9,9,7,9,9,9,9,9,10,9,9,9,8,9,11,9,7,9,9,10,9,9,9,8,8,9,9,9,9,9,9,7,8,8,9,4,4,6)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
